# Fruta / hortaliza, semanal
# Insert a new week's worth of data (4 quality grades: Especial, Primera,
# Segunda, Tercera) at row 552, pushing the existing rows 552:563 down to
# 556:567.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the current row 552, shifting everything below
# (rows 552:563) down to 556:567.
$ws.Range("A552:T555").Insert()

# Values that are constant across this whole data block.
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100101
$producto   = "Berries"
$categoriaId = 100112025
$categoria   = "Frutilla"
$variedad    = "Sin especificar"
$unidad      = "`$/bandeja 7 kilos"
$origen      = "Provincia de Melipilla"
$kgUnidad    = 7

# New rows of data for the 2021-09-09 (serial 44448) report.
$newRows = @(
    @{ Row = 552; Calidad = "Especial"; Volumen = 150; PMin = 25000; PMax = 25000; PProm = 25000; PKg = 3571 },
    @{ Row = 553; Calidad = "Primera";  Volumen = 200; PMin = 20000; PMax = 20000; PProm = 20000; PKg = 2857 },
    @{ Row = 554; Calidad = "Segunda";  Volumen = 120; PMin = 15000; PMax = 15000; PProm = 15000; PKg = 2143 },
    @{ Row = 555; Calidad = "Tercera";  Volumen = 100; PMin = 10000; PMax = 10000; PProm = 10000; PKg = 1429 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = 44448
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
